$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This worksheet holds a rolling window of gyroscope samples (x, y, z) in
# columns A:C, rows 2:21 (row 1 is the header "x","y","z"). A new capture
# of 4 samples was pushed into the top of the window, shifting the
# existing samples down by 4 rows and dropping the 4 oldest samples off
# the bottom so the window keeps its fixed size (A1:C21).

# Final (post-shift) values for every data row, top to bottom: the first
# 4 rows are the newly captured samples, the rest are the previous
# samples shifted down by 4 rows (the last 4 original rows fall off the
# end of the window and are discarded).
$data = @(
    @(0.03629761248826986, 0.01907121278345579, 0.05546045627444995),
    @(0.05165476366877556, -0.0003787364251911958, 0.03377473920583711),
    @(0.02702467799186697, -0.02729956846684218, -0.006963863894343374),
    @(0.02345722466707222, -0.01078177168965329, -0.01979203335940831),
    @(0.02312735825777049, -0.06875288158655167, -0.001032362207770397),
    @(-0.1697350136935719, -0.2189157873392111, -0.04578435219824339),
    @(-0.2920058012008669, -0.2727635514736176, 0.126033713221551),
    @(-0.2655186891555784, -0.2190562760829923, 0.6186354464292545),
    @(-0.2926044583320622, -0.3912286460399645, 0.472412636876103),
    @(-0.02654820919036799, 0.1870408368110673, -0.7209736722707767),
    @(1.565036740303041, -0.7936239337921156, -2.525304698944093),
    @(1.211272468566887, -0.4785640525817794, -2.630715656280514),
    @(0.08948565244675255, -0.522057590484637, -1.111322727203364),
    @(0.5279402807354899, -1.442680406570422, 0.08100073337555302),
    @(-0.09629679918289313, 0.4833410131931261, 0.534000061750409),
    @(-0.2217929553985592, -0.1976026952266694, -0.07598552823066797),
    @(-0.01921781659126266, -0.1601200026273725, -0.3930673503875731),
    @(0.1233581319451331, 0.0740124344825745, -0.1553430366516109),
    @(0.02880229651927913, 0.08823948577046314, 0.1134376771748063),
    @(-0.04516126751899735, -0.04690834142267682, 0.04470311729237415)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
